# Fruta / hortaliza, semanal
# Insert two new weekly observations (a new "Especial" and "Primera" quality
# record for the most recent Monday) at the top of the data block that starts
# at row 479, pushing all the older rows down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 479.
$ws.Rows.Item(479).Insert()
$ws.Rows.Item(479).Insert()

# ---- Row 479: new "Especial" quality record ----
$ws.Cells.Item(479, 1).Value = 4
$ws.Cells.Item(479, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(479, 3).Value = "Los Lagos"
$ws.Cells.Item(479, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(479, 4).Value = 44939
$ws.Cells.Item(479, 5).Value = 10
$ws.Cells.Item(479, 6).Value = "Fruta"
$ws.Cells.Item(479, 7).Value = 100102
$ws.Cells.Item(479, 8).Value = "Cítricos"
$ws.Cells.Item(479, 9).Value = 100102006
$ws.Cells.Item(479, 10).Value = "Pomelo"
$ws.Cells.Item(479, 11).Value = "Start Ruby"
$ws.Cells.Item(479, 12).Value = "Especial"
$ws.Cells.Item(479, 13).Value = 100
$ws.Cells.Item(479, 14).Value = 15000
$ws.Cells.Item(479, 15).Value = 15000
$ws.Cells.Item(479, 16).Value = 15000
$ws.Cells.Item(479, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(479, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(479, 19).Value = 1071
$ws.Cells.Item(479, 20).Value = 14

# ---- Row 480: new "Primera" quality record ----
$ws.Cells.Item(480, 1).Value = 4
$ws.Cells.Item(480, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(480, 3).Value = "Los Lagos"
$ws.Cells.Item(480, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(480, 4).Value = 44939
$ws.Cells.Item(480, 5).Value = 10
$ws.Cells.Item(480, 6).Value = "Fruta"
$ws.Cells.Item(480, 7).Value = 100102
$ws.Cells.Item(480, 8).Value = "Cítricos"
$ws.Cells.Item(480, 9).Value = 100102006
$ws.Cells.Item(480, 10).Value = "Pomelo"
$ws.Cells.Item(480, 11).Value = "Start Ruby"
$ws.Cells.Item(480, 12).Value = "Primera"
$ws.Cells.Item(480, 13).Value = 200
$ws.Cells.Item(480, 14).Value = 12000
$ws.Cells.Item(480, 15).Value = 13000
$ws.Cells.Item(480, 16).Value = 12500
$ws.Cells.Item(480, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(480, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(480, 19).Value = 893
$ws.Cells.Item(480, 20).Value = 14
